$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1), new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match formatting of the existing header cells (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-10
$values = @(
    @(8, 8),
    @(8, 9),
    @(8, 9),
    @(8, 8),
    @(6, 6),
    @(6, 6),
    @(8, 8),
    @(6, 6),
    @(3, 3)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $values[$i][0]
    $ws.Cells.Item($row, 10).Value = $values[$i][1]
}
